$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "72.325.46"
Set-TextValue "E2" "  +0.12%  "
Set-TextValue "D3" "2.653.40"
Set-TextValue "E3" "  +0.46%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "592.26"
Set-TextValue "D6" "174.81"
Set-TextValue "E6" "  -3.01%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "E8" "  -0.73%  "
Set-TextValue "D9" "2.651.46"
Set-TextValue "E9" "  +0.42%  "
Set-TextValue "D10" "0.171"
Set-TextValue "E10" "  -2.18%  "
Set-TextValue "E11" "  +1.43%  "
Set-TextValue "E12" "  -0.33%  "
Set-TextValue "E13" "  -1.57%  "
Set-TextValue "D14" "3.141.55"
Set-TextValue "E14" "  +0.95%  "
Set-TextValue "D15" "0.0000186"
Set-TextValue "E15" "  -2.12%  "
Set-TextValue "D16" "72.254.49"
Set-TextValue "E16" "  +0.15%  "
Set-TextValue "D17" "26.05"
Set-TextValue "D18" "2.609.79"
Set-TextValue "E18" "  -1.34%  "
Set-TextValue "D19" "12.33"
Set-TextValue "E19" "  +3.02%  "
Set-TextValue "D20" "7.99"
Set-TextValue "E20" "  +0.43%  "
Set-TextValue "D21" "372.20"
Set-TextValue "E21" "  -1.71%  "
Set-TextValue "D22" "4.17"
Set-TextValue "E22" "  -0.56%  "
Set-TextValue "D23" "2.07"
Set-TextValue "E23" "  -0.82%  "
Set-TextValue "D24" "71.70"
Set-TextValue "E24" "  -2.29%  "
Set-TextValue "E25" "  -0.16%  "
Set-TextValue "D26" "4.27"
Set-TextValue "E26" "  -3.26%  "
Set-TextValue "D27" "9.68"
Set-TextValue "E27" "  -4.23%  "
Set-TextValue "D28" "2.793.71"
Set-TextValue "E28" "  +0.39%  "
Set-TextValue "E29" "  -0.19%  "
Set-TextValue "D30" "0.0₃0957"
Set-TextValue "E30" "  +0.01%  "
Set-TextValue "D31" "8.09"
Set-TextValue "E31" "  -0.66%  "
Set-TextValue "D32" "499.07"
Set-TextValue "E32" "  -4.89%  "
Set-TextValue "E33" "  -2.67%  "
Set-TextValue "E34" "  -0.89%  "
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "D36" "161.88"
Set-TextValue "E37" "  +3.96%  "
Set-TextValue "D38" "19.40"
Set-TextValue "E38" "  +0.15%  "
Set-TextValue "E39" "  -0.95%  "
Set-TextValue "E40" "  -3.26%  "
Set-TextValue "E41" "  +0.01%  "
Set-TextValue "E42" "  -6.41%  "
Set-TextValue "D43" "2.56"
Set-TextValue "E43" "  -3.05%  "
Set-TextValue "D44" "4.90"
Set-TextValue "E44" "  -3.80%  "
Set-TextValue "E45" "  -1.43%  "
Set-TextValue "D46" "39.16"
Set-TextValue "E46" "  -0.71%  "
Set-TextValue "D47" "154.27"
Set-TextValue "E47" "  +1.93%  "
Set-TextValue "E48" "  -0.96%  "
Set-TextValue "E49" "  +0.92%  "
Set-TextValue "D50" "1.69"
Set-TextValue "E50" "  -0.44%  "
Set-TextValue "D51" "0.0748"
Set-TextValue "E51" "  -1.52%  "
